$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (List exp.: NR/GCS/C/O)
$ws.Range("B2").Value = 2.82352941176471
$ws.Range("C2").Value = 2.13612565445026
$ws.Range("D2").Value = 2.67361111111111
$ws.Range("E2").Value = 2.585

# Row 3 (List exp.: NR/C/O)
$ws.Range("B3").Value = 2.06081081081081
$ws.Range("C3").Value = 1.58415841584158
$ws.Range("D3").Value = 1.94074074074074
$ws.Range("E3").Value = 2.02702702702703
$ws.Range("F3").Value = 1.38845553822153

# Row 4 (List exp.: GCS/C/O)
$ws.Range("B4").Value = 2.05142857142857
$ws.Range("C4").Value = 1.71929824561404
$ws.Range("D4").Value = 1.97857142857143
$ws.Range("E4").Value = 1.82706766917293
$ws.Range("F4").Value = 1.40233236151604

# Row 5 (List exp.: C/O)
$ws.Range("B5").Value = 1.31055900621118
$ws.Range("C5").Value = 0.927461139896373
$ws.Range("D5").Value = 1.04838709677419
$ws.Range("E5").Value = 1.31901840490798
$ws.Range("F5").Value = 0.894428152492669
